$wb = $excel.ActiveWorkbook

function Format-Header($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

# --- Sheet "save" ---
$ws1 = $wb.Worksheets.Item("save")
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "month_id"
$ws1.Range("C1").Value = "week_id"
$ws1.Range("D1").Value = "day_id"
$ws1.Range("E1").Value = "day"
$ws1.Range("F1").Value = "sum"
$ws1.Range("G1").Value = "date"
Format-Header $ws1.Range("A1:G1")

# --- Sheet "save_log" ---
$ws2 = $wb.Worksheets.Item("save_log")
$ws2.Range("A1").Value = "save_id"
$ws2.Range("B1").Value = "desc"
$ws2.Range("C1").Value = "amount"
$ws2.Range("D1").Value = "time"
Format-Header $ws2.Range("A1:D1")
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 54
$ps2.RightMargin = 54
$ps2.TopMargin = 72
$ps2.BottomMargin = 72
$ps2.HeaderMargin = 36
$ps2.FooterMargin = 36

# --- Sheet "cost" ---
$ws3 = $wb.Worksheets.Item("cost")
$ws3.Range("A1").Value = "id"
$ws3.Range("B1").Value = "month_id"
$ws3.Range("C1").Value = "week_id"
$ws3.Range("D1").Value = "day_id"
$ws3.Range("E1").Value = "day"
$ws3.Range("F1").Value = "sum"
$ws3.Range("G1").Value = "date"
Format-Header $ws3.Range("A1:G1")
$ps3 = $ws3.PageSetup
$ps3.LeftMargin = 54
$ps3.RightMargin = 54
$ps3.TopMargin = 72
$ps3.BottomMargin = 72
$ps3.HeaderMargin = 36
$ps3.FooterMargin = 36

# --- Sheet "cost_log" ---
$ws4 = $wb.Worksheets.Item("cost_log")
$ws4.Range("A1").Value = "cost_id"
$ws4.Range("B1").Value = "desc"
$ws4.Range("C1").Value = "amount"
$ws4.Range("D1").Value = "time"
Format-Header $ws4.Range("A1:D1")
$ps4 = $ws4.PageSetup
$ps4.LeftMargin = 54
$ps4.RightMargin = 54
$ps4.TopMargin = 72
$ps4.BottomMargin = 72
$ps4.HeaderMargin = 36
$ps4.FooterMargin = 36

# --- Selections per sheet ---
[void]$ws1.Range("D13").Select()
[void]$ws2.Range("B29").Select()
[void]$ws3.Range("F13").Select()
[void]$ws4.Range("C5").Select()

# --- Active sheet / tab ---
[void]$ws1.Activate()
[void]$ws1.Range("D13").Select()
